$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.595.68"
$ws.Range("E2").Value = "  +1.98%  "

$ws.Range("D3").Value = "3.932.05"
$ws.Range("E3").Value = "  +0.52%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'535.66"
$ws.Range("E5").Value = "  +9.79%  "

$ws.Range("D6").Value = "'145.20"
$ws.Range("E6").Value = "  -1.14%  "

$ws.Range("D7").Value = "'0.619"
$ws.Range("E7").Value = "  -0.66%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").Value = "'0.728"
$ws.Range("E9").Value = "  -0.08%  "

$ws.Range("E10").Value = "  +3.91%  "

$ws.Range("D11").Value = "'0.0000337"
$ws.Range("E11").Value = "  -2.03%  "

$ws.Range("D12").Value = "'42.62"
$ws.Range("E12").Value = "  -1.03%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'10.39"
$ws.Range("E13").Value = "  -4.92%  "

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "4.560.11"
$ws.Range("E14").Value = "  +0.54%  "

$ws.Range("D15").Value = "3.936.48"
$ws.Range("E15").Value = "  -2.04%  "

$ws.Range("E16").Value = "  -1.57%  "

$ws.Range("E17").Value = "  -0.23%  "

$ws.Range("E18").Value = "  +7.29%  "

$ws.Range("D19").Value = "'19.81"
$ws.Range("E19").Value = "  -0.57%  "

$ws.Range("D20").Value = "69.464.23"
$ws.Range("E20").Value = "  +1.62%  "

$ws.Range("D21").Value = "'432.34"
$ws.Range("E21").Value = "  -0.28%  "

$ws.Range("D22").Value = "'3.38"
$ws.Range("E22").Value = "  -4.37%  "

$ws.Range("B23").Value = "PancakeSwap"
$ws.Range("C23").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D23").Value = "'4.23"
$ws.Range("E23").Value = "  +17.11%  "

$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").Value = "'14.51"
$ws.Range("E24").Value = "  -2.33%  "

$ws.Range("D25").Value = "'88.65"
$ws.Range("E25").Value = "  +0.91%  "

$ws.Range("D26").Value = "'11.81"
$ws.Range("E26").Value = "  +4.15%  "

$ws.Range("D27").Value = "'10.80"
$ws.Range("E27").Value = "  -4.07%  "

$ws.Range("D28").Value = "'36.62"
$ws.Range("E28").Value = "  -3.83%  "

$ws.Range("D29").Value = "'700.69"
$ws.Range("E29").Value = "  -2.90%  "

$ws.Range("D30").Value = "'13.32"
$ws.Range("E30").Value = "  -3.21%  "

$ws.Range("D31").Value = "'72.39"
$ws.Range("E31").Value = "  +19.96%  "

$ws.Range("E32").Value = "  -1.99%  "

$ws.Range("D34").Value = "'0.465"
$ws.Range("E34").Value = "  +16.79%  "

$ws.Range("D35").Value = "'6.12"
$ws.Range("E35").Value = "  -2.32%  "

$ws.Range("B36").Value = "InjectiveProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D36").Value = "'40.52"
$ws.Range("E36").Value = "  -2.58%  "

$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").Value = "0.0₃0858"
$ws.Range("E37").Value = "  -1.31%  "

$ws.Range("E38").Value = "  +0.74%  "

$ws.Range("E39").Value = "  +0.12%  "

$ws.Range("E40").Value = "  -0.07%  "

$ws.Range("D41").Value = "'0.0482"
$ws.Range("E41").Value = "  -0.06%  "

$ws.Range("D42").Value = "'2.82"
$ws.Range("E42").Value = "  -3.60%  "

$ws.Range("E43").Value = "  +6.72%  "

$ws.Range("E44").Value = "  -4.69%  "

$ws.Range("D45").Value = "'3.18"
$ws.Range("E45").Value = "  +12.60%  "

$ws.Range("D46").Value = "'0.142"
$ws.Range("E46").Value = "  +0.70%  "

$ws.Range("D47").Value = "'3.36"
$ws.Range("E47").Value = "  +0.57%  "

$ws.Range("D48").Value = "0.0₆0351"
$ws.Range("E48").Value = "  -1.03%  "

$ws.Range("E49").Value = "  -2.02%  "

$ws.Range("D50").Value = "'144.93"
$ws.Range("E50").Value = "  +0.29%  "

$ws.Range("E51").Value = "  -2.05%  "

